# college-list.xlsx: append 16 new college rows (554-569) to the dashboard,
# extend the UPPER() helper formula down to match, and recolor the new
# "city" helper column (E) entries with the dark-grey (#222222) font used
# for the newest batch of entries ("color change in dashboard").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows of data -------------------------------------------------
# Columns: A = college name, B = UPPER(A) helper formula,
#          C = sequence number, E = city, F = state

$rows = @(
    @{ R = 554; A = 'GD Rungta College of Engineering and Technology (GDRCET) , Bhilai'; C = 557; E = 'Bhilai'; F = 'Chattisgarh' },
    @{ R = 555; A = 'Bharath university , Chennai'; C = 558; E = 'Chennai'; F = 'Tamil Nadu' },
    @{ R = 556; A = 'Thapar University , Patiala'; C = 559; E = 'Patiala'; F = 'Punjab' },
    @{ R = 557; A = 'Jai Prakash University (JP), Chapra'; C = 560; E = 'Chapra'; F = 'Bihar' },
    @{ R = 558; A = 'Assam Don Bosco University - Azara Campus , Guwahati'; C = 561; E = 'Guwahati'; F = 'Assam' },
    @{ R = 559; A = 'Babulal Tarabai Institute of Research and Technology (BTIRT), Sagar'; C = 562; E = 'Sagar'; F = 'Madhya Pradesh' },
    @{ R = 560; A = 'Anand Engineering College (AEC) , Agra'; C = 563; E = 'Agra'; F = 'Uttar Pradesh' },
    @{ R = 561; A = 'PDM College of Engineering (PDMCE) , Bahadurgarh'; C = 564; E = 'Bahadurgarh'; F = 'Haryana' },
    @{ R = 562; A = 'Government Model Engineering College (MEC) , Kochi'; C = 565; E = 'Kochin'; F = 'Kerala' },
    @{ R = 563; A = 'International Institute of Information Technology (IIIT), Bhubaneswar'; C = 566; E = 'Bhubaneswar'; F = 'Orissa' },
    @{ R = 564; A = 'Ideal Institute of Technology , Ghaziabad'; C = 567; E = 'Ghaziabad'; F = 'Uttar Pradesh' },
    @{ R = 565; A = 'Delhi Institute of Tool Engineering , Delhi'; C = 568; E = 'Delhi'; F = 'Delhi' },
    @{ R = 566; A = 'National Institute of Technology (NIT) , Delhi'; C = 569; E = 'Delhi'; F = 'Delhi' },
    @{ R = 567; A = 'IGNOU , MUZZAFARPUR'; C = 570; E = 'Muzaffarpur'; F = 'Bihar' },
    @{ R = 568; A = 'IGNOU , Delhi'; C = 571; E = 'Delhi'; F = 'Delhi' },
    @{ R = 569; A = 'MBICEM GGSIPU , DELHI'; C = 572; E = 'Delhi'; F = 'Delhi' }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
}

# --- Extend the shared UPPER() helper formula down column B -----------
$ws.Range("B554:B569").Formula = "=UPPER(A554)"

# --- Recolor the new "city" helper cells (dark grey, #222222) ---------
# This is the visible "color change in dashboard" referenced in the commit.
$ws.Range("E561").Font.Color = 2236962

# --- Update the active selection to reflect where editing stopped -----
$ws.Range("A570").Select()
